# Weekly update: insert a new data row at row 24 (pushing the existing
# row 24 record down to row 25), then populate the new row 24 with the
# latest weekly price record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 24; this shifts the current row 24
# (and everything below it) down by one row, so the old row 24 becomes
# row 25, preserving its formatting/values automatically.
$ws.Rows("24:24").Insert()

# Populate the newly inserted row 24 with the new weekly record.
$ws.Cells.Item(24, 1).Value = 7
$ws.Cells.Item(24, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(24, 3).Value = "Ñuble"

$ws.Cells.Item(24, 4).Value = 44610
$ws.Cells.Item(24, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Cells.Item(24, 5).Value = 16
$ws.Cells.Item(24, 6).Value = 100112040
$ws.Cells.Item(24, 7).Value = "Cilantro"
$ws.Cells.Item(24, 8).Value = "Sin especificar"
$ws.Cells.Item(24, 9).Value = "Primera"
$ws.Cells.Item(24, 10).Value = 100
$ws.Cells.Item(24, 11).Value = 550
$ws.Cells.Item(24, 12).Value = 600
$ws.Cells.Item(24, 13).Value = 575
$ws.Cells.Item(24, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(24, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(24, 16).Value = 575
$ws.Cells.Item(24, 17).Value = 1
$ws.Cells.Item(24, 18).Value = "Hortaliza"
